$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
Write-Host $ws.Name
Write-Host $wb.Worksheets.Count
for ($i=1; $i -le $wb.Worksheets.Count; $i++) {
    Write-Host $wb.Worksheets.Item($i).Name
}
